# Refresh the cryptocurrency price/volume table (Price = column D, Volume(1h) = column E).
# Some "Price" values (e.g. "0.566") look like plain numbers to Excel, which would silently
# convert them from text to a Number cell on assignment. To keep them as text (matching the
# source data, which stores every price/volume cell as a string), those cells are briefly
# switched to a Text number format before the value is written, then ClearFormats() restores
# the cell's original (default/general) formatting so no stray number format lingers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.752.77"
$ws.Range("D3").Value = "2.541.20"
$ws.Range("E3").Value = "  -0.14%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "311.38"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.47%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "100.05"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +0.76%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.566"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -0.73%  "
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.522"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -1.81%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.38"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -1.44%  "
$ws.Range("E11").Value = "  -0.25%  "
$ws.Range("E12").Value = "  +0.97%  "
$ws.Range("E13").Value = "  -1.43%  "
$ws.Range("D14").Value = "2.933.34"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.39"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -3.56%  "
$ws.Range("D16").Value = "2.544.27"
$ws.Range("E16").Value = "  +0.59%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.814"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -3.03%  "
$ws.Range("D18").Value = "42.770.22"
$ws.Range("E18").Value = "  +0.06%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.72"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -0.26%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.32"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -0.52%  "
$ws.Range("D21").Value = "0.0₃0950"
$ws.Range("E21").Value = "  -0.62%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "69.75"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +0.44%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "242.92"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -1.95%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.87"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -1.30%  "
$ws.Range("E25").Value = "  -1.74%  "
$ws.Range("E26").Value = "  +0.09%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "25.69"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -3.29%  "
$ws.Range("E28").Value = "  -1.05%  "
$ws.Range("E29").Value = "  +0.18%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "38.23"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -4.32%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.90"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +2.89%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "157.68"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -0.83%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.71"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +3.48%  "
$ws.Range("E34").Value = "  +1.86%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0791"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -1.48%  "
$ws.Range("E36").Value = "  -4.13%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "17.86"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -2.79%  "
$ws.Range("E38").Value = "  -5.31%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.110"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -0.82%  "
$ws.Range("E40").Value = "  -0.28%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.12"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -0.62%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "21.81"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -2.85%  "
$ws.Range("E43").Value = "  +0.15%  "
$ws.Range("E44").Value = "  -0.30%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.25"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +0.57%  "
$ws.Range("D46").Value = "1.995.39"
$ws.Range("E46").Value = "  +0.35%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.08"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +0.44%  "
$ws.Range("D48").Value = "2.785.74"
$ws.Range("E49").Value = "  -0.97%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "79.63"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -2.06%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "72.35"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -1.46%  "
